$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1397.7646
$ws.Range("I98").Value = 950.5714
$ws.Range("K98").Value = 950.5714
$ws.Range("M98").Value = 547.4286
$ws.Range("H99").Value = 316.5
$ws.Range("I99").Value = 316.5
$ws.Range("K99").Value = 949.5
$ws.Range("M99").Value = 548.5
$ws.Range("H113").Value = 5058
$ws.Range("J113").Value = 5116.5
$ws.Range("L113").Value = 5116.5
$ws.Range("N113").Value = -11624.5
$ws.Range("H122").Value = 1397.7646
$ws.Range("I122").Value = 950.5714
$ws.Range("K122").Value = 2851.7142
$ws.Range("M122").Value = -401.7142000000003
$ws.Range("H132").Value = 1093.4615
$ws.Range("I132").Value = 992.5
$ws.Range("J132").Value = 1430
$ws.Range("K132").Value = 2977.5
$ws.Range("L132").Value = 4290
$ws.Range("M132").Value = -447.5
$ws.Range("N132").Value = -9350
$ws.Range("H137").Value = 3825.2666
$ws.Range("I137").Value = 2488.9
$ws.Range("K137").Value = 7466.700000000001
$ws.Range("M137").Value = -4916.700000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 5000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -4827
$ws.Range("N6").ClearContents()
$ws.Range("H32").Value = 5009.946
$ws.Range("I32").Value = 4968.758
$ws.Range("K32").Value = 4968.758
$ws.Range("M32").Value = -4681.758
$ws.Range("H61").Value = 1175.5
$ws.Range("I61").Value = 1203.6923
$ws.Range("K61").Value = 1203.6923
$ws.Range("M61").Value = -991.6922999999999
$ws.Range("H74").Value = 2232.5312
$ws.Range("I74").Value = 2278.1
$ws.Range("K74").Value = 2278.1
$ws.Range("M74").Value = -1404.1
$ws.Range("H77").Value = 2232.5312
$ws.Range("I77").Value = 2278.1
$ws.Range("K77").Value = 11390.5
$ws.Range("M77").Value = -7022.5
$ws.Range("H132").Value = 3150.9
$ws.Range("I132").Value = 2287.2856
$ws.Range("J132").Value = 5166
$ws.Range("K132").Value = 6861.8568
$ws.Range("L132").Value = 15498
$ws.Range("M132").Value = -4331.8568
$ws.Range("N132").Value = -20558
$ws.Range("H136").Value = 1175.5
$ws.Range("I136").Value = 1203.6923
$ws.Range("K136").Value = 3611.0769
$ws.Range("M136").Value = -1061.0769

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2704.85
$ws.Range("I134").Value = 2784.0527
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 8352.158100000001
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -5817.158100000001
$ws.Range("N134").Value = -8670
$ws.Range("H135").Value = 53836.418
$ws.Range("J135").Value = 53836.418
$ws.Range("L135").Value = 53836.418
$ws.Range("N135").Value = -63976.418

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3065.25
$ws.Range("I31").Value = 2998.375
$ws.Range("J31").Value = 3199
$ws.Range("K31").Value = 2998.375
$ws.Range("L31").Value = 3199
$ws.Range("M31").Value = -2703.375
$ws.Range("N31").Value = -3789
$ws.Range("H34").Value = 3065.25
$ws.Range("I34").Value = 2998.375
$ws.Range("J34").Value = 3199
$ws.Range("K34").Value = 2998.375
$ws.Range("L34").Value = 3199
$ws.Range("M34").Value = -2796.375
$ws.Range("N34").Value = -3603
$ws.Range("H58").Value = 1854.2727
$ws.Range("J58").Value = 3014
$ws.Range("L58").Value = 3014
$ws.Range("N58").Value = -3420
$ws.Range("H136").Value = 1854.2727
$ws.Range("J136").Value = 3014
$ws.Range("L136").Value = 9042
$ws.Range("N136").Value = -14142

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 200000
$ws.Range("J37").Value = 200000
$ws.Range("L37").Value = 600000
$ws.Range("N37").Value = -600224
$ws.Range("H64").Value = 1000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 1000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H80").Value = 501943.5
$ws.Range("I80").Value = 3888
$ws.Range("K80").Value = 11664
$ws.Range("M80").Value = -10728
$ws.Range("H83").Value = 501943.5
$ws.Range("I83").Value = 3888
$ws.Range("K83").Value = 34992
$ws.Range("M83").Value = -30312
$ws.Range("H113").Value = 1245.375
$ws.Range("I113").Value = 1448
$ws.Range("J113").Value = 1177.8334
$ws.Range("K113").Value = 4344
$ws.Range("L113").Value = 3533.5002
$ws.Range("M113").Value = -2174
$ws.Range("N113").Value = -7873.5002
$ws.Range("H136").Value = 4120
$ws.Range("I136").Value = 4120
$ws.Range("K136").Value = 12360
$ws.Range("M136").Value = -7260
$ws.Range("H137").Value = 5149.5557
$ws.Range("I137").Value = 11633.333
$ws.Range("J137").Value = 1907.6666
$ws.Range("K137").Value = 34899.999
$ws.Range("L137").Value = 5722.9998
$ws.Range("M137").Value = -29799.999
$ws.Range("N137").Value = -15922.9998
$ws.Range("H140").Value = 2236.125
$ws.Range("I140").Value = 2236.125
$ws.Range("K140").Value = 6708.375
$ws.Range("M140").Value = -1528.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1419.8695
$ws.Range("I132").Value = 729.1875
$ws.Range("J132").Value = 2998.5715
$ws.Range("K132").Value = 2187.5625
$ws.Range("L132").Value = 8995.7145
$ws.Range("M132").Value = 342.4375
$ws.Range("N132").Value = -14055.7145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3768
$ws.Range("I40").Value = 3768
$ws.Range("K40").Value = 3768
$ws.Range("M40").Value = -3632
$ws.Range("H122").Value = 2965
$ws.Range("I122").Value = 2965
$ws.Range("K122").Value = 8895
$ws.Range("M122").Value = -6445
$ws.Range("H132").Value = 3206.6538
$ws.Range("I132").Value = 2171.2727
$ws.Range("K132").Value = 6513.8181
$ws.Range("M132").Value = -3983.8181
$ws.Range("H136").Value = 3340.682
$ws.Range("I136").Value = 3251.8948
$ws.Range("J136").Value = 3903
$ws.Range("K136").Value = 9755.6844
$ws.Range("L136").Value = 11709
$ws.Range("M136").Value = -7205.6844
$ws.Range("N136").Value = -16809

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2294.4048
$ws.Range("I132").Value = 1543.069
$ws.Range("K132").Value = 4629.207
$ws.Range("M132").Value = -2099.207
$ws.Range("H136").Value = 1232.3793
$ws.Range("I136").Value = 913.0769
$ws.Range("J136").Value = 3999.6667
$ws.Range("K136").Value = 2739.2307
$ws.Range("L136").Value = 11999.0001
$ws.Range("M136").Value = -189.2307000000001
$ws.Range("N136").Value = -17099.0001
